$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.17"
$ws.Range("E2").Value = "'-0.78%"
$ws.Range("D3").Value = "'31.45"
$ws.Range("E3").Value = "'-1.97%"
$ws.Range("D4").Value = "'5.118"
$ws.Range("E4").Value = "'-2.26%"
$ws.Range("D5").Value = "'0.07352"
$ws.Range("E5").Value = "'-2.63%"
$ws.Range("D6").Value = "'2.209"
$ws.Range("E6").Value = "'44.75%"
$ws.Range("D7").Value = "'7.924"
$ws.Range("E7").Value = "'-0.27%"
$ws.Range("D8").Value = "'3.794"
$ws.Range("E8").Value = "'-0.68%"
$ws.Range("D9").Value = "'0.9197"
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("D10").Value = "'0.1706"
$ws.Range("E10").Value = "'0.72%"
$ws.Range("D11").Value = "'0.07509"
$ws.Range("E11").Value = "'-4.99%"
$ws.Range("D12").Value = "'0.08155"
$ws.Range("E12").Value = "'1.75%"
$ws.Range("D13").Value = "'0.03026"
$ws.Range("E13").Value = "'-0.34%"
$ws.Range("D14").Value = "'0.09947"
$ws.Range("E14").Value = "'0.44%"
$ws.Range("D15").Value = "'0.001512"
$ws.Range("E15").Value = "'1.20%"
$ws.Range("D16").Value = "'0.006144"
$ws.Range("E16").Value = "'-5.75%"
$ws.Range("D17").Value = "'3.452"
$ws.Range("E17").Value = "'0.16%"
$ws.Range("D18").Value = "'2.222"
$ws.Range("E18").Value = "'-0.29%"
$ws.Range("D19").Value = "'0.3279"
$ws.Range("E19").Value = "'-0.64%"
$ws.Range("D20").Value = "'0.1338"
$ws.Range("E20").Value = "'0.78%"
$ws.Range("D21").Value = "'4.664"
$ws.Range("E21").Value = "'4.33%"
$ws.Range("D22").Value = "'0.04639"
$ws.Range("E22").Value = "'1.02%"
$ws.Range("D23").Value = "'0.1569"
$ws.Range("E23").Value = "'-3.06%"
$ws.Range("D24").Value = "'0.001229"
$ws.Range("E24").Value = "'1.12%"
$ws.Range("D25").Value = "'0.004468"
$ws.Range("E25").Value = "'-0.05%"
$ws.Range("D26").Value = "'0.0001304"
$ws.Range("E26").Value = "'-6.74%"
$ws.Range("D27").Value = "'0.0003437"
$ws.Range("E27").Value = "'92.90%"
$ws.Range("D39").Value = "'0.01728"
$ws.Range("E39").Value = "'2.01%"
$ws.Range("D40").Value = "'0.04517"
$ws.Range("E40").Value = "'0.96%"
$ws.Range("D41").Value = "'0.007269"
$ws.Range("E41").Value = "'4.35%"
$ws.Range("E42").Value = "'-0.59%"
$ws.Range("D43").Value = "'0.002236"
$ws.Range("E43").Value = "'7.67%"
$ws.Range("D44").Value = "'0.01066"
$ws.Range("E44").Value = "'-22.28%"
$ws.Range("D45").Value = "'0.00006293"
$ws.Range("E45").Value = "'2.10%"

# Row 46: becomes BOLO data
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D46").Value = "'0.8085"
$ws.Range("E46").Value = "'-56.17%"

# Row 47: becomes CoinbaseStockToken data
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D47").Value = "'0.009992"
$ws.Range("E47").Value = "'-22.99%"
